# Updated cryptos list on Thu Aug 15 20:07:34 UTC 2024 with GitHub Actions
#
# Refreshes the Price (D) / Volume(1h) (E) columns for each coin row, and
# swaps the Monero / EthereumClassic rows (33/34) and the RenderToken /
# VeChain rows (50/51) to reflect the new ranking order (Coin/Link/Price
# swap together; the rank number in column A is unaffected).
#
# For D-column prices that are valid numeric literals (e.g. "512.43"),
# Excel's Range.Value setter auto-coerces the string to a real number,
# which both changes the stored cell type and reformats the text (loses
# trailing zeros, switches to float precision). The source data stores
# these as literal text, so for those cells we force text storage via
# NumberFormat "@" before assigning, then restore the cell's default
# ("Normal") style afterwards so no stray number formatting is left
# behind - only the underlying value type (text) sticks.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.025.26'
$ws.Range('E2').Value = '  -3.27%  '

$ws.Range('D3').Value = '2.540.63'
$ws.Range('E3').Value = '  -4.64%  '

$ws.Range('E4').Value = '  +0.04%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '512.43'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -2.38%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '139.72'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.21%  '

$ws.Range('E7').Value = '  +0.05%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.557'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -2.19%  '

$ws.Range('E9').Value = '  -7.72%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0989'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -3.55%  '

$ws.Range('E11').Value = '  -3.67%  '

$ws.Range('E12').Value = '  -0.12%  '

$ws.Range('D13').Value = '2.993.90'
$ws.Range('E13').Value = '  -4.34%  '

$ws.Range('D14').Value = '57.024.83'
$ws.Range('E14').Value = '  -3.28%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '20.01'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -5.13%  '

$ws.Range('E16').Value = '  -3.10%  '

$ws.Range('D17').Value = '2.532.94'
$ws.Range('E17').Value = '  -5.91%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '332.04'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.95%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '4.27'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -2.45%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '10.06'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -3.14%  '

$ws.Range('E21').Value = '  -4.12%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '0.999'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.18%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '64.04'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.30%  '

$ws.Range('E24').Value = '  -0.55%  '

$ws.Range('E25').Value = '  +0.01%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.399'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -4.49%  '

$ws.Range('D27').Value = '2.655.90'
$ws.Range('E27').Value = '  -4.32%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.90'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.66%  '

$ws.Range('D29').Value = '0.0₃0749'
$ws.Range('E29').Value = '  -6.38%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.01%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '6.26'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -6.31%  '

$ws.Range('E32').Value = '  -2.90%  '

$ws.Range('B33').Value = 'Monero'
$ws.Range('C33').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '148.34'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -1.69%  '

$ws.Range('B34').Value = 'EthereumClassic'
$ws.Range('C34').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '18.43'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  -2.24%  '

$ws.Range('E35').Value = '  -4.89%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '1.12'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -5.68%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.839'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -5.56%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '35.77'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  -2.89%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.820'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -5.66%  '

$ws.Range('E40').Value = '  -3.04%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '1.00'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.12%  '

$ws.Range('E42').Value = '  -3.77%  '

$ws.Range('E43').Value = '  -1.36%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '10.62'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.39%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.577'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -6.47%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '258.04'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -6.30%  '

$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '0.0518'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.48%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '18.42'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -7.63%  '

$ws.Range('D49').Value = '1.973.36'
$ws.Range('E49').Value = '  -3.75%  '

$ws.Range('B50').Value = 'RenderToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '4.52'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -3.88%  '

$ws.Range('B51').Value = 'VeChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0220'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  -3.43%  '

